$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1 + 2: split "...demandé pour se choix" into two runs with the
# "_GoBack" bookmark sitting right before "se choix". Word keeps a single,
# unique "_GoBack" bookmark that auto-relocates to the last edit point, so
# adding it here also removes it from its old spot (before "Lancer
# l'installateur").
# ---------------------------------------------------------------------------
$rGoBack = $d.Content
$rGoBack.Find.Execute("se choix", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$goBackPoint = $d.Range($rGoBack.Start, $rGoBack.Start)
$d.Bookmarks.Add("_GoBack", $goBackPoint) | Out-Null

# ---------------------------------------------------------------------------
# Change 3: shrink " à partir du raccourci sur le bureau, du raccourci dans
# le menu démarrer ou à partir de" down to " à partir de", while keeping the
# text split into the two runs " à partir" and " de" (matching the target
# run layout, separate from the untouched "Lancer l'exécutable" run and the
# untouched italic path run that follows).
# ---------------------------------------------------------------------------

# 1) Delete the middle portion that isn't wanted any more.
$rMid = $d.Content
$middle = " du raccourci sur le bureau, du raccourci dans le menu démarrer ou à partir"
$rMid.Find.Execute($middle, $false, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null

# 2) Re-split the now-merged run right after "Lancer l'exécutable" so that
#    " à partir de" becomes its own run again.
$rAfterExec = $d.Content
$rAfterExec.Find.Execute("Lancer l’exécutable", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$splitAfterExec = $d.Range($rAfterExec.End, $rAfterExec.End)
$d.Bookmarks.Add("ZZZ_TMP_SPLIT_1", $splitAfterExec) | Out-Null
$d.Bookmarks("ZZZ_TMP_SPLIT_1").Delete() | Out-Null

# 3) Split " à partir de" into " à partir" and " de".
$rPartir = $d.Content
$rPartir.Find.Execute(" à partir", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$splitPartir = $d.Range($rPartir.End, $rPartir.End)
$d.Bookmarks.Add("ZZZ_TMP_SPLIT_2", $splitPartir) | Out-Null
$d.Bookmarks("ZZZ_TMP_SPLIT_2").Delete() | Out-Null

Write-Output "done"
